$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handback-status entry: 5c7432d3-9c05-4341-b2f1-924c95c9f677.md
# Adds a 4th row of data to each of the 3 sheets (Overview, zh-cn, de-de)
# ---------------------------------------------------------------------------

$fileGuid    = "5c7432d3-9c05-4341-b2f1-924c95c9f677"
$fileName    = "$fileGuid.md"
$pathName    = "e2e\$fileGuid.md"
$pathNameUrl = "e2e/$fileGuid.md"
$statusText  = "Handed back: in sync with en-US"

# =================== Sheet "Overview" (table3 / A1:G) ======================
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A4").Value = $fileName
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = $statusText
$wsOv.Range("F4").Value = $statusText
$wsOv.Range("G4").Value = "2016-10-13 13:57:18"
$wsOv.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add(
    $wsOv.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af6da57df78594503dbd058d30d799a650731141/$pathNameUrl",
    "",
    "",
    $pathName
) | Out-Null

# =================== Sheet "zh-cn" (table1 / A1:P) ==========================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$zhXlf = "$fileGuid.28806e771721b21f9aed8add520f4e9bd3749c9c.zh-cn.xlf"

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusText
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = "2016-10-13 13:57:07"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = "2016-10-13 13:57:50"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af6da57df78594503dbd058d30d799a650731141/$pathNameUrl",
    "",
    "",
    $fileName
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/28806e771721b21f9aed8add520f4e9bd3749c9c/$pathNameUrl",
    "",
    "",
    $fileName
) | Out-Null

# =================== Sheet "de-de" (table2 / A1:P) ==========================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$deXlf = "$fileGuid.28806e771721b21f9aed8add520f4e9bd3749c9c.de-de.xlf"

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusText
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = "2016-10-13 13:57:18"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = "2016-10-13 13:58:06"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af6da57df78594503dbd058d30d799a650731141/$pathNameUrl",
    "",
    "",
    $fileName
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/28806e771721b21f9aed8add520f4e9bd3749c9c/$pathNameUrl",
    "",
    "",
    $fileName
) | Out-Null
